$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.132.50'
$ws.Range('E2').Value = '  +0.49%  '
$ws.Range('D3').Value = '1.655.40'
$ws.Range('E3').Value = '  +0.19%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.93'
$ws.Range('E5').Value = '  +0.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5296'
$ws.Range('E6').Value = '  +1.79%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2614'
$ws.Range('E8').Value = '  -0.36%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06328'
$ws.Range('E9').Value = '  +1.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.43'
$ws.Range('E10').Value = '  -0.69%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07803'
$ws.Range('E11').Value = '  +1.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.518'
$ws.Range('E12').Value = '  +1.06%  '
$ws.Range('D13').Value = '1.654.43'
$ws.Range('E13').Value = '  -0.65%  '
$ws.Range('D14').Value = '1.883.08'
$ws.Range('E14').Value = '  +0.24%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5489'
$ws.Range('E15').Value = '  +1.15%  '
$ws.Range('D16').Value = '0.0₅8212'
$ws.Range('E16').Value = '  +1.49%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.41'
$ws.Range('E17').Value = '  +1.07%  '
$ws.Range('D18').Value = '26.125.98'
$ws.Range('E18').Value = '  +0.40%  '
$ws.Range('E19').Value = '  -0.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.600'
$ws.Range('E20').Value = '  +0.61%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '191.12'
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('E22').Value = '  +0.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.028'
$ws.Range('E23').Value = '  +0.81%  '
$ws.Range('E24').Value = '  +0.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.02'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1227'
$ws.Range('E26').Value = '  -0.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.217'
$ws.Range('E27').Value = '  -0.35%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.00'
$ws.Range('E28').Value = '  -0.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.455'
$ws.Range('E29').Value = '  +3.86%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05773'
$ws.Range('E30').Value = '  -2.76%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.273'
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.559'
$ws.Range('E32').Value = '  +1.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.271'
$ws.Range('E33').Value = '  +0.70%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.605'
$ws.Range('E34').Value = '  +3.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.801'
$ws.Range('E35').Value = '  +1.62%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9512'
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.413'
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5764'
$ws.Range('E38').Value = '  +1.50%  '
$ws.Range('E39').Value = '  +1.07%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8560'
$ws.Range('E40').Value = '  +0.99%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.784'
$ws.Range('E41').Value = '  -1.68%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '104.70'
$ws.Range('E42').Value = '  +4.10%  '
$ws.Range('D43').Value = '1.045.63'
$ws.Range('E43').Value = '  +4.67%  '
$ws.Range('E44').Value = '  +0.14%  '
$ws.Range('D45').Value = '1.797.67'
$ws.Range('E45').Value = '  +0.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.90'
$ws.Range('E46').Value = '  +0.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.003'
$ws.Range('E47').Value = '  -0.52%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4340'
$ws.Range('E48').Value = '  +0.87%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.895'
$ws.Range('E49').Value = '  -1.02%  '
$ws.Range('E50').Value = '  +0.03%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.445'
$ws.Range('E51').Value = '  -1.59%  '
